# Fig5c_PW_PB_bar.xlsx update
# Correct the misspelled "Anticipitory" condition label to "Anticipatory"
# (the data itself - all numeric values and the other row/column labels -
# is unchanged; only the spelling of this one label is fixed) and restore
# the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Anticipatory"

$ws.Range("A3").Select()
